# Weekly fruit/vegetable data update: insert a new weekly price record for
# "Espinaca" (Vega Central Mapocho de Santiago) as row 431, pushing all the
# following rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 431 (shifts rows 431..455 down to 432..456,
# and carries the column D date-number-format down with them).
$ws.Rows.Item(431).Insert()

# Populate the new row with the new weekly record. All of the "fixed"
# columns (market/region/category/unit/origin/classification) match the
# rest of this sheet's rows; only the date + quality/volume/price columns
# are specific to this new entry.
$ws.Cells.Item(431, 1).Value = 9
$ws.Cells.Item(431, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(431, 3).Value = "Metropolitana"
$ws.Cells.Item(431, 4).Value = 44753
$ws.Cells.Item(431, 5).Value = 13
$ws.Cells.Item(431, 6).Value = 100112012
$ws.Cells.Item(431, 7).Value = "Espinaca"
$ws.Cells.Item(431, 8).Value = "Sin especificar"
$ws.Cells.Item(431, 9).Value = "Primera"
$ws.Cells.Item(431, 10).Value = 70
$ws.Cells.Item(431, 11).Value = 11000
$ws.Cells.Item(431, 12).Value = 11000
$ws.Cells.Item(431, 13).Value = 11000
$ws.Cells.Item(431, 14).Value = "$/cuna 10 kilos"
$ws.Cells.Item(431, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(431, 16).Value = 1100
$ws.Cells.Item(431, 17).Value = 10
$ws.Cells.Item(431, 18).Value = "Hortaliza"
